{"js": "// Replace the page title (appears twice: the H1 heading and the bold\n// \"title\" run near the end of the document).\nconst titleResults = context.document.body.search(\n  \"Play G.R.O.W Slot for Free \\u2013 Review of Gameplay Mechanics\",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < titleResults.items.length; i++) {\n  titleResults.items[i].insertText(\n    \"Play G.R.O.W Slot Game for Free\",\n    Word.InsertLocation.replace\n  );\n}\n\n// Simple one-to-one text replacements (search text -> new text).\nconst replacements = [\n  [\"Fun and unique gameplay mechanics\", \"Free to play\"],\n  [\"Green-themed graphics and animations\", \"Multiple bonus features\"],\n  [\"Four bonus features that trigger frequently\", \"Mobile compatibility\"],\n  [\"Moderately profitable RTP of 96.47%\", \"Appealing graphics and animation\"],\n  [\"No progressive jackpot feature\", \"May not appeal to experienced slot players\"],\n  [\n    \"Find out how to play G.R.O.W online slot for free. Read our unbiased review of the gameplay mechanics, RTP, and bonus features of this unique garden-themed slot game.\",\n    \"Read our review of G.R.O.W, a free online slot game with multiple bonus features.\",\n  ],\n];\n\nconst searchResultsList = [];\nfor (const [oldText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  searchResultsList.push(results);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const results = searchResultsList[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# wdReplaceAll so both occurrences of the page title (the H1 heading and\n# the bold \"title\" run near the end of the document) get updated.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n\n$enDash = [char]0x2013\n\nReplace-Text (\"Play G.R.O.W Slot for Free \" + $enDash + \" Review of Gameplay Mechanics\") \"Play G.R.O.W Slot Game for Free\"\nReplace-Text \"Fun and unique gameplay mechanics\" \"Free to play\"\nReplace-Text \"Green-themed graphics and animations\" \"Multiple bonus features\"\nReplace-Text \"Four bonus features that trigger frequently\" \"Mobile compatibility\"\nReplace-Text \"Moderately profitable RTP of 96.47%\" \"Appealing graphics and animation\"\nReplace-Text \"No progressive jackpot feature\" \"May not appeal to experienced slot players\"\nReplace-Text \"Find out how to play G.R.O.W online slot for free. Read our unbiased review of the gameplay mechanics, RTP, and bonus features of this unique garden-themed slot game.\" \"Read our review of G.R.O.W, a free online slot game with multiple bonus features.\"\n"}
